$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsDR = $wb.Worksheets.Item("DR")

# --- "About" sheet updates ---

$wsAbout.Activate()

# B3 Source text value stays the same content, only shared-string index shuffled (no value change)
# Row 9: "Notes:" header stays same text.

# Replace old "Annual Perc" note (row 10) with the new multi-line explanatory note (rows 10-14)
$wsAbout.Range("A10").Value = "This is the annual percentage rate by which future savings (e.g. fuel cost savings) are discounted when"
$wsAbout.Range("A11").Value = "making price-driven purchasing decisions in the current year. The value used should be one that is"
$wsAbout.Range("A12").Value = "reasonable for people who are looking to buy fuel-consuming capital equipment, such as industrial"
$wsAbout.Range("A13").Value = "equipment or building components. The model works in real dollars, so this rate should be the growth"
$wsAbout.Range("A14").Value = "in real value, not the growth in nominal value plus real value."

# Row 15 stays blank, then rows 16-18 carry the notes that used to be rows 10-11
$wsAbout.Range("A16").Value = "We choose to use a 3% discount rate here, for consistency with the 3% rate used for the central estimate"
$wsAbout.Range("A17").Value = "of Social Cost of Carbon (in the SCoC variable), as well as the discount rate built into the health"
$wsAbout.Range("A18").Value = "damages values in the SCoHIbP Social Cost of Health Impacts by Pollutant variable."

# --- "DR" sheet updates ---

$wsDR.Activate()

$wsDR.Range("B1").Value = "Annual Perc (dimensionless)"
$wsDR.Range("B1").WrapText = $true
$wsDR.Rows("1:1").RowHeight = 30

$wsDR.Range("B1").Select() | Out-Null

# --- Restore "About" as the active sheet with its final selection ---

$wsAbout.Activate()
$wsAbout.Range("A16:A18").Select() | Out-Null
